# The commit removes one post entry from the sheet:
#   "「僕たちは協力して不可能を達成するよ」" (originally row 694).
# Deleting that entire row naturally shifts every following row up by
# one position (695 -> 694, 696 -> 695, ... 816 -> 815) and shrinks the
# used range from A1:C816 down to A1:C815, exactly matching the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(694).Delete()
